$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell B2: replace multi-line placeholder with single summary line ---
$ws.Range("B2").Value = "Component: Multi-Functional Tool Application"

# --- Test case rows 6-30 (columns B-H) ---
# Row 6
$ws.Cells.Item(6,2).Value = "TC001"
$ws.Cells.Item(6,3).Value = "Application is installed on a desktop computer"
$ws.Cells.Item(6,4).Value = "Verify application installation on desktop"
$ws.Cells.Item(6,5).Value = "1. Copy MultiFunctionalTool_For_Desktop.zip from specified tec-share location`n2. Extract the contents`n3. Double-click on MultiFunctionalToolApplication"
$ws.Cells.Item(6,6).Value = "Application launches successfully"
$ws.Cells.Item(6,7).Value = ""
$ws.Cells.Item(6,8).Value = ""

# Row 7
$ws.Cells.Item(7,2).Value = "TC002"
$ws.Cells.Item(7,3).Value = "Application is installed on a laptop"
$ws.Cells.Item(7,4).Value = "Verify application installation on laptop"
$ws.Cells.Item(7,5).Value = "1. Copy MultiFunctionalTool_For_Laptop.zip from specified tec-share location`n2. Extract the contents`n3. Double-click on MultiFunctionalToolApplication"
$ws.Cells.Item(7,6).Value = "Application launches successfully"
$ws.Cells.Item(7,7).Value = ""
$ws.Cells.Item(7,8).Value = ""

# Row 8
$ws.Cells.Item(8,2).Value = "TC003"
$ws.Cells.Item(8,3).Value = "Application is installed and running"
$ws.Cells.Item(8,4).Value = "Verify Network Packet Capture start functionality"
$ws.Cells.Item(8,5).Value = "1. Navigate to Network Packet Capture section`n2. Click on Start button"
$ws.Cells.Item(8,6).Value = "Packet capture begins successfully"
$ws.Cells.Item(8,7).Value = ""
$ws.Cells.Item(8,8).Value = ""

# Row 9
$ws.Cells.Item(9,2).Value = "TC004"
$ws.Cells.Item(9,3).Value = "Network Packet Capture is running"
$ws.Cells.Item(9,4).Value = "Verify Network Packet Capture stop functionality"
$ws.Cells.Item(9,5).Value = "1. Navigate to Network Packet Capture section`n2. Click on Stop button"
$ws.Cells.Item(9,6).Value = "1. Packet capture stops`n2. .pcap file is generated`n3. File is copied to MFP's Shared Folder`n4. Shared Folder opens automatically"
$ws.Cells.Item(9,7).Value = ""
$ws.Cells.Item(9,8).Value = ""

# Row 10
$ws.Cells.Item(10,2).Value = "TC005"
$ws.Cells.Item(10,3).Value = "Application is installed and running"
$ws.Cells.Item(10,4).Value = "Verify Memory Leak Check functionality"
$ws.Cells.Item(10,5).Value = "1. Navigate to Memory Leak Check section`n2. Select a protocol`n3. Run the memory leak check"
$ws.Cells.Item(10,6).Value = "Memory Leak Comparison Table is displayed with accurate information"
$ws.Cells.Item(10,7).Value = ""
$ws.Cells.Item(10,8).Value = ""

# Row 11
$ws.Cells.Item(11,2).Value = "TC006"
$ws.Cells.Item(11,3).Value = "Application is installed and running"
$ws.Cells.Item(11,4).Value = "Verify Debug Log Collection functionality"
$ws.Cells.Item(11,5).Value = "1. Navigate to Debug Log Collection section`n2. Click on Run button"
$ws.Cells.Item(11,6).Value = "1. Script executes successfully`n2. Logs are collected`n3. Logs are copied to MFP's Shared Folder`n4. Shared Folder opens automatically"
$ws.Cells.Item(11,7).Value = ""
$ws.Cells.Item(11,8).Value = ""

# Row 12
$ws.Cells.Item(12,2).Value = "TC007"
$ws.Cells.Item(12,3).Value = "Debug Log Collection has been run once with empty folder"
$ws.Cells.Item(12,4).Value = "Verify Debug Log Collection retry functionality"
$ws.Cells.Item(12,5).Value = "1. Navigate to Debug Log Collection section`n2. Click on Run button again"
$ws.Cells.Item(12,6).Value = "Logs are successfully collected and displayed in the Shared Folder"
$ws.Cells.Item(12,7).Value = ""
$ws.Cells.Item(12,8).Value = ""

# Row 13
$ws.Cells.Item(13,2).Value = "TC008"
$ws.Cells.Item(13,3).Value = "Application is installed and running"
$ws.Cells.Item(13,4).Value = "Verify Diagnostic Code Details for ECC"
$ws.Cells.Item(13,5).Value = "1. Navigate to Diagnostic Code Details section`n2. Select ECC option"
$ws.Cells.Item(13,6).Value = "Relevant ECC diagnostic details are displayed"
$ws.Cells.Item(13,7).Value = ""
$ws.Cells.Item(13,8).Value = ""

# Row 14
$ws.Cells.Item(14,2).Value = "TC009"
$ws.Cells.Item(14,3).Value = "Application is installed and running"
$ws.Cells.Item(14,4).Value = "Verify Diagnostic Code Details for Network Protocols"
$ws.Cells.Item(14,5).Value = "1. Navigate to Diagnostic Code Details section`n2. Select Network Protocols option"
$ws.Cells.Item(14,6).Value = "Relevant Network Protocols diagnostic details are displayed"
$ws.Cells.Item(14,7).Value = ""
$ws.Cells.Item(14,8).Value = ""

# Row 15
$ws.Cells.Item(15,2).Value = "TC010"
$ws.Cells.Item(15,3).Value = "Application is installed and running"
$ws.Cells.Item(15,4).Value = "Verify Diagnostic Code Details for High Security Mode"
$ws.Cells.Item(15,5).Value = "1. Navigate to Diagnostic Code Details section`n2. Select High Security Mode option"
$ws.Cells.Item(15,6).Value = "Relevant High Security Mode diagnostic details are displayed"
$ws.Cells.Item(15,7).Value = ""
$ws.Cells.Item(15,8).Value = ""

# Row 16
$ws.Cells.Item(16,2).Value = "TC011"
$ws.Cells.Item(16,3).Value = "Application is installed and running"
$ws.Cells.Item(16,4).Value = "Verify Diagnostic Code Details for commonly used codes"
$ws.Cells.Item(16,5).Value = "1. Navigate to Diagnostic Code Details section`n2. Select a commonly used diagnostic code"
$ws.Cells.Item(16,6).Value = "Relevant job-specific details are displayed"
$ws.Cells.Item(16,7).Value = ""
$ws.Cells.Item(16,8).Value = ""

# Row 17
$ws.Cells.Item(17,2).Value = "TC012"
$ws.Cells.Item(17,3).Value = "Application is installed and running"
$ws.Cells.Item(17,4).Value = "Verify 08 Diagnostic Code Value Get functionality"
$ws.Cells.Item(17,5).Value = "1. Navigate to 08 Diagnostic Code Value section`n2. Select a diagnostic code`n3. Click on Get button"
$ws.Cells.Item(17,6).Value = "Current value of the selected diagnostic code is displayed"
$ws.Cells.Item(17,7).Value = ""
$ws.Cells.Item(17,8).Value = ""

# Row 18
$ws.Cells.Item(18,2).Value = "TC013"
$ws.Cells.Item(18,3).Value = "Application is installed and running"
$ws.Cells.Item(18,4).Value = "Verify 08 Diagnostic Code Value Set functionality"
$ws.Cells.Item(18,5).Value = "1. Navigate to 08 Diagnostic Code Value section`n2. Select a diagnostic code`n3. Enter a new value`n4. Click on Set button"
$ws.Cells.Item(18,6).Value = "Diagnostic code value is updated successfully"
$ws.Cells.Item(18,7).Value = ""
$ws.Cells.Item(18,8).Value = ""

# Row 19
$ws.Cells.Item(19,2).Value = "TC014"
$ws.Cells.Item(19,3).Value = "Application is installed and running"
$ws.Cells.Item(19,4).Value = "Verify Protocol Configuration Get functionality"
$ws.Cells.Item(19,5).Value = "1. Navigate to Protocol Configuration section`n2. Open Protocol Selection Window`n3. Select a protocol`n4. Click on Get button"
$ws.Cells.Item(19,6).Value = "Current protocol value is displayed"
$ws.Cells.Item(19,7).Value = ""
$ws.Cells.Item(19,8).Value = ""

# Row 20
$ws.Cells.Item(20,2).Value = "TC015"
$ws.Cells.Item(20,3).Value = "Application is installed and running"
$ws.Cells.Item(20,4).Value = "Verify Protocol Configuration Set functionality"
$ws.Cells.Item(20,5).Value = "1. Navigate to Protocol Configuration section`n2. Open Protocol Selection Window`n3. Select a protocol`n4. Enter a new value`n5. Click on Set button"
$ws.Cells.Item(20,6).Value = "Message indicating `"Set protocol values operation still needs to be implemented`" is displayed"
$ws.Cells.Item(20,7).Value = ""
$ws.Cells.Item(20,8).Value = "Feature not yet implemented"

# Row 21
$ws.Cells.Item(21,2).Value = "TC016"
$ws.Cells.Item(21,3).Value = "Application is running with multiple features active"
$ws.Cells.Item(21,4).Value = "Verify application performance under load"
$ws.Cells.Item(21,5).Value = "1. Start Network Packet Capture`n2. Run Memory Leak Check`n3. Collect Debug Logs simultaneously"
$ws.Cells.Item(21,6).Value = "All operations complete successfully without significant delay or application crash"
$ws.Cells.Item(21,7).Value = ""
$ws.Cells.Item(21,8).Value = ""

# Row 22
$ws.Cells.Item(22,2).Value = "TC017"
$ws.Cells.Item(22,3).Value = "Application is installed and running"
$ws.Cells.Item(22,4).Value = "Verify GUI responsiveness"
$ws.Cells.Item(22,5).Value = "1. Navigate through all sections of the application`n2. Click on various buttons and options"
$ws.Cells.Item(22,6).Value = "GUI responds promptly to all user interactions with no visible lag"
$ws.Cells.Item(22,7).Value = ""
$ws.Cells.Item(22,8).Value = ""
$ws.Rows.Item(22).AutoFit()

# Row 23
$ws.Cells.Item(23,2).Value = "TC018"
$ws.Cells.Item(23,3).Value = "Application is installed on a system with minimum specifications"
$ws.Cells.Item(23,4).Value = "Verify application performance on low-end systems"
$ws.Cells.Item(23,5).Value = "1. Launch the application`n2. Test all major functionalities"
$ws.Cells.Item(23,6).Value = "Application runs without significant performance issues"
$ws.Rows.Item(23).AutoFit()

# Row 24
$ws.Cells.Item(24,2).Value = "TC019"
$ws.Cells.Item(24,3).Value = "Application is running"
$ws.Cells.Item(24,4).Value = "Verify error handling for network disconnection"
$ws.Cells.Item(24,5).Value = "1. Start a network-dependent operation`n2. Disconnect from network during operation"
$ws.Cells.Item(24,6).Value = "Application displays appropriate error message and handles the disconnection gracefully"
$ws.Rows.Item(24).AutoFit()

# Row 25
$ws.Cells.Item(25,2).Value = "TC020"
$ws.Cells.Item(25,3).Value = "Application is running"
$ws.Cells.Item(25,4).Value = "Verify compatibility with different Linux systems"
$ws.Cells.Item(25,5).Value = "1. Install application on different Linux distributions`n2. Test all major functionalities"
$ws.Cells.Item(25,6).Value = "Application works consistently across different Linux environments"
$ws.Rows.Item(25).AutoFit()

# Row 26
$ws.Cells.Item(26,2).Value = "TC021"
$ws.Cells.Item(26,3).Value = "Application is running"
$ws.Cells.Item(26,4).Value = "Verify usability for first-time users"
$ws.Cells.Item(26,5).Value = "1. Have a first-time user navigate through the application`n2. Ask them to perform basic tasks without instructions"
$ws.Cells.Item(26,6).Value = "User can successfully navigate and perform basic tasks with minimal confusion"
$ws.Rows.Item(26).AutoFit()

# Row 27
$ws.Cells.Item(27,2).Value = "TC022"
$ws.Cells.Item(27,3).Value = "Application is running"
$ws.Cells.Item(27,4).Value = "Verify application behavior during long operations"
$ws.Cells.Item(27,5).Value = "1. Start a time-consuming operation (e.g., extensive packet capture)`n2. Monitor application behavior"
$ws.Cells.Item(27,6).Value = "Application remains responsive and provides progress indication"
$ws.Rows.Item(27).AutoFit()

# Row 28
$ws.Cells.Item(28,2).Value = "TC023"
$ws.Cells.Item(28,3).Value = "Application is running"
$ws.Cells.Item(28,4).Value = "Verify data integrity of collected logs"
$ws.Cells.Item(28,5).Value = "1. Collect debug logs`n2. Verify the content and format of collected logs"
$ws.Cells.Item(28,6).Value = "Logs are complete, properly formatted, and contain all required information"
$ws.Rows.Item(28).AutoFit()

# Row 29
$ws.Cells.Item(29,2).Value = "TC024"
$ws.Cells.Item(29,3).Value = "Application is running"
$ws.Cells.Item(29,4).Value = "Verify memory usage during extended operation"
$ws.Cells.Item(29,5).Value = "1. Run the application for an extended period (8+ hours)`n2. Monitor memory usage"
$ws.Cells.Item(29,6).Value = "Application does not exhibit memory leaks or excessive memory consumption"
$ws.Rows.Item(29).AutoFit()

# Row 30
$ws.Cells.Item(30,2).Value = "TC025"
$ws.Cells.Item(30,3).Value = "Application is running"
$ws.Cells.Item(30,4).Value = "Verify application recovery after crash"
$ws.Cells.Item(30,5).Value = "1. Force the application to crash`n2. Restart the application"
$ws.Cells.Item(30,6).Value = "Application restarts properly and recovers previous state if applicable"
$ws.Rows.Item(30).AutoFit()

# --- Dimension: the authoring process touched column L at some point, which expands
# the recorded sheet dimension from A1:I31 to A1:L31 even though no visible data lives there ---
$ws.Range("L31").Value = "x"
$ws.Range("L31").ClearContents()
